$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must remain plain text (they use "." as thousands
# separators or would otherwise be re-interpreted as numbers by Excel),
# so force a text number format before assigning, then restore the default style.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.446.61"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.329.10"
$ws.Range("D3").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.24"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.43"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("D7").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.905.08"
$ws.Range("D12").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.453.09"
$ws.Range("D14").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.327.39"
$ws.Range("D16").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "431.75"
$ws.Range("D18").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.22"
$ws.Range("D22").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.70"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.455.88"
$ws.Range("D25").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.198"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.04"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.22"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.63"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.83"
$ws.Range("D37").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.890.47"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.73"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.768"
$ws.Range("D42").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.35"
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.01"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.32"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.50"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "317.93"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0272"
$ws.Range("D50").Style = "Normal"

# Coin name / link / volume columns assign safely as text already.
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +7.96%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("E51").Value = "  +4.74%  "
